$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-9: "Förändrad" date (column C) moves forward one day, 45183 -> 45184
foreach ($r in 2..9) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# Row 2 hyperlink formulas gain a second (display-text) argument.
# S2 keeps the broken/unbalanced quoting exactly as produced by the
# source update (the trailing quote from the URL literal is not closed,
# and the added label text is appended raw).
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_BRACKE/artfynd/A 30834-2023.xlsx, "A 30834-2023"")'

$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_BRACKE/kartor/A 30834-2023.png", "A 30834-2023")'

$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_BRACKE/klagomål/A 30834-2023.docx", "A 30834-2023")'

$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_BRACKE/klagomålsmail/A 30834-2023.docx", "A 30834-2023")'

$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_BRACKE/tillsyn/A 30834-2023.docx", "A 30834-2023")'

# Y2 was stored as literal inline text (not a real formula); turn it into
# an actual formula cell with the same two-argument HYPERLINK call.
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_BRACKE/tillsynsmail/A 30834-2023.docx", "A 30834-2023")'
